$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (Email), shifting Email..Grant Date right by one.
# This mirrors the user inserting a "Department" column in Excel.
$ws.Columns.Item(3).Insert()

# The hyperlinks that lived on the old column C (Email) do not automatically
# follow the shift performed by Columns.Insert() in this environment, so
# capture them first and re-create them on the new column D.
$hlAddresses = @()
$hlDisplays = @()
foreach ($hl in $ws.Hyperlinks) {
    $hlAddresses += $hl.Address
    $hlDisplays += $hl.TextToDisplay
}
$ws.Hyperlinks.Delete()

# Populate the new "Department" column header and values.
$ws.Range("C1").Value = "Department"
$ws.Range("C2").Value = "Tech"
$ws.Range("C3").Value = "HR"
$ws.Range("C4").Value = "Ops"
$ws.Range("C5").Value = "Tech"
$ws.Range("C6").Value = "HR"
$ws.Range("C7").Value = "Ops"

# Re-create the hyperlinks on column D (now the Email column) pointing at the
# same addresses/display text as before, restoring the original (non-default)
# hyperlink look: Arial 11, blue, no underline - matching the rest of the sheet
# instead of Excel's auto-applied "Hyperlink" theme style.
for ($i = 0; $i -lt $hlAddresses.Count; $i++) {
    $row = $i + 2
    $rng = $ws.Range("D" + $row)
    $ws.Hyperlinks.Add($rng, $hlAddresses[$i], "", "", $hlDisplays[$i])
    $rng.Font.Name = "Arial"
    $rng.Font.Size = 11
    $rng.Font.Color = 16711680
    $rng.Font.Underline = $false
}

# Move the active selection the way the source workbook ended up (C8).
$ws.Range("C8").Select()
